# Apply cell-level updates per the diff (values scraped from coinranking.com).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.311.02"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "1.733.15"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.66"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.11"
$ws.Range("E8").Value = "  +10.33%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "1.741.01"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.70"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "28.277.93"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.48"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "0.0₃0754"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.66"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +3.59%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.27"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.489.94"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.66"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.603"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.46"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "1.882.50"
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.800"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.74"
$ws.Range("E47").Value = "  +9.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.04"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +7.99%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.20"
$ws.Range("E51").Value = "  +0.66%  "
